# Add a new worksheet "JUANCITO" after the existing "PINI" sheet,
# with header row Producto / Codigo in A1:B1.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "JUANCITO"

$newSheet.Range("A1").Value = "Producto"
$newSheet.Range("B1").Value = "Codigo"
